# Sum_r.xlsx edit: add a third sheet ("Sheet3") with a new data series
# (random sphere packing bed, MR ...) and extend Sheet2 with two more
# result columns (I, J) plus fill in the previously-header-only H column.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# 1) Sheet2: fill in column H (header already present, data was missing)
#    and add two brand-new columns I and J, with matching headers.
# ---------------------------------------------------------------------
$ws2.Range("H1").NumberFormat = $ws2.Range("G1").NumberFormat
$ws2.Range("I1").NumberFormat = $ws2.Range("G1").NumberFormat
$ws2.Range("J1").NumberFormat = $ws2.Range("G1").NumberFormat

$ws2.Range("I1").Value = "Ibm-mrta"
$ws2.Range("J1").Value = "Psm-mrtb"

$ws2.Range("H2:J7").NumberFormat = $ws2.Range("G2").NumberFormat

$ws2.Range("H2").Value = 0.129487491
$ws2.Range("I2").Value = 0.123421751
$ws2.Range("J2").Value = 0.0914734671

$ws2.Range("H3").Value = 0.0871822565
$ws2.Range("I3").Value = 0.0888339192
$ws2.Range("J3").Value = 0.0713024869

$ws2.Range("H4").Value = 0.0822136402
$ws2.Range("I4").Value = 0.0827071371
$ws2.Range("J4").Value = 0.0708284843

$ws2.Range("H5").Value = 0.0799363019
$ws2.Range("I5").Value = 0.0800679722
$ws2.Range("J5").Value = 0.0709763945

$ws2.Range("H6").Value = 0.0785600114
$ws2.Range("I6").Value = 0.0785638005
$ws2.Range("J6").Value = 0.0713362252

$ws2.Range("H7").Value = 0.0775220646
$ws2.Range("I7").Value = 0.0774896726
$ws2.Range("J7").Value = 0.0714295318

# Move the selection/cursor off the old A26 -> A2 on Sheet2 before the
# new sheet becomes active (matches the committed selection state).
[void]$ws2.Range("A2").Select()

# ---------------------------------------------------------------------
# 2) Add "Sheet3" after Sheet2 (copy Sheet2 so the page setup / header /
#    footer / column-width metadata carries over), then strip the
#    copied B:J data and rebuild the small 2-column data set.
# ---------------------------------------------------------------------
$ws2.Copy([System.Type]::Missing, $ws2)
$ws3 = $wb.Worksheets.Item("Sheet2 (2)")
$ws3.Name = "Sheet3"

$ws3.Range("B1:J7").Clear()

$ws3.Range("G1").Value = "Psm-mrta"
$ws3.Range("H1").Value = "Ibm-mrtb"
$ws3.Range("I1").Value = "Ibm-mrta"
$ws3.Range("J1").Value = "Psm-mrtb"

$ws3.Range("A1").Value = 1.6
$ws3.Range("A2").Value = 10
$ws3.Range("A3").Value = 30
$ws3.Range("A4").Value = 50
$ws3.Range("A5").Value = 70
$ws3.Range("A6").Value = 90
$ws3.Range("A7").Value = 110

[void]$ws3.Range("A2").Select()
